$d = $word.ActiveDocument

# Locate the paragraph that holds "ABlasterCharacter" (end of the
# Blaster/ACharacter/ABlasterCharacter list branch).
$anchor = $null
$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "ABlasterCharacter") {
        $anchor = $p
        $anchorIndex = $i
    }
    $i = $i + 1
}

# Insert a new paragraph right after it; it inherits the anchor's
# paragraph/run formatting (pStyle, numPr/numId, spacing, rPr lang, etc.)
$anchor.Range.InsertParagraphAfter()

# ListLevelNumber is 1-based (w:ilvl = ListLevelNumber - 1), so ilvl=2 -> 3.
$controllerPara = $d.Paragraphs.Item($anchorIndex + 2)
$controllerPara.Range.ListFormat.ListLevelNumber = 3
$controllerPara.Range.Text = "APlayerController"

# Insert the next paragraph after that one, at ilvl=3 -> ListLevelNumber=4.
$controllerPara.Range.InsertParagraphAfter()
$blasterControllerPara = $d.Paragraphs.Item($anchorIndex + 3)
$blasterControllerPara.Range.ListFormat.ListLevelNumber = 4
$blasterControllerPara.Range.Text = "ABlasterPlayerController"
